$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4.917497143191863
$ws.Range("D2").Value = 8.37953447089845
$ws.Range("E2").Value = 13.1715048820428
$ws.Range("F2").Value = 38.31587132213537
$ws.Range("G2").Value = 3.637461528759324
$ws.Range("J2").Value = 9.970737365885174
$ws.Range("N2").Value = 18.99769123883691
$ws.Range("O2").Value = 29.15539173534993
$ws.Range("C3").Value = 4.750038287576572
$ws.Range("D3").Value = 8.394650609564671
$ws.Range("E3").Value = 13.15978120347579
$ws.Range("F3").Value = 37.81025459191146
$ws.Range("G3").Value = 3.641781056497854
$ws.Range("J3").Value = 9.975187776357476
$ws.Range("N3").Value = 18.40031508502705
$ws.Range("O3").Value = 28.79417684634628
$ws.Range("C4").Value = 4.64591576783808
$ws.Range("D4").Value = 8.404845278276248
$ws.Range("E4").Value = 13.15545742062762
$ws.Range("F4").Value = 37.50866947976713
$ws.Range("G4").Value = 3.644570570457484
$ws.Range("J4").Value = 9.979858037231841
$ws.Range("N4").Value = 18.02485520896364
$ws.Range("O4").Value = 28.57929408398432
$ws.Range("C5").Value = 4.603233074408617
$ws.Range("D5").Value = 8.40922952907969
$ws.Range("E5").Value = 13.15441859798503
$ws.Range("F5").Value = 37.3881454899031
$ws.Range("G5").Value = 3.645741988542309
$ws.Range("J5").Value = 9.982247850571035
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 28.49356298653657
$ws.Range("C6").Value = 4.596132744288094
$ws.Range("D6").Value = 8.409971417741623
$ws.Range("E6").Value = 13.1542897787859
$ws.Range("F6").Value = 37.36828010956902
$ws.Range("G6").Value = 3.645938599671039
$ws.Range("J6").Value = 9.982674048260407
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 28.47944106643438
$ws.Range("C7").Value = 4.645341050300152
$ws.Range("D7").Value = 8.404903474949895
$ws.Range("E7").Value = 13.15544048259847
$ws.Range("F7").Value = 37.50703424964986
$ws.Range("G7").Value = 3.644586228050326
$ws.Range("J7").Value = 9.979888297599556
$ws.Range("N7").Value = 18.02277304766463
$ws.Range("O7").Value = 28.57813032523592
$ws.Range("C8").Value = 4.860071148188108
$ws.Range("D8").Value = 8.384557071697326
$ws.Range("E8").Value = 13.16686625684063
$ws.Range("F8").Value = 38.13978326231944
$ws.Range("G8").Value = 3.638922489429729
$ws.Range("J8").Value = 9.971869385334447
$ws.Range("N8").Value = 18.79364780656866
$ws.Range("O8").Value = 29.02947253667056
$ws.Range("C9").Value = 5.267797834203673
$ws.Range("D9").Value = 8.351896488627711
$ws.Range("E9").Value = 13.21204874189638
$ws.Range("F9").Value = 39.44433865126174
$ws.Range("G9").Value = 3.628898923999232
$ws.Range("J9").Value = 9.971540063062042
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 29.96479493695725
$ws.Range("C10").Value = 5.555577424452109
$ws.Range("D10").Value = 8.332303123471856
$ws.Range("E10").Value = 13.25905958193109
$ws.Range("F10").Value = 40.43258108286911
$ws.Range("G10").Value = 3.622185812989481
$ws.Range("J10").Value = 9.980703773332708
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 30.67632919590986
$ws.Range("C11").Value = 5.683289480366905
$ws.Range("D11").Value = 8.324343274162553
$ws.Range("E11").Value = 13.28342094219014
$ws.Range("F11").Value = 40.88667420338377
$ws.Range("G11").Value = 3.619271314870455
$ws.Range("J11").Value = 9.986915426135772
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 31.00394998365633
$ws.Range("C12").Value = 5.731144146241925
$ws.Range("D12").Value = 8.321465978278923
$ws.Range("E12").Value = 13.29307089710891
$ws.Range("F12").Value = 41.05911770167053
$ws.Range("G12").Value = 3.618187556162777
$ws.Range("J12").Value = 9.989561056539859
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 31.12846399086969
$ws.Range("C13").Value = 5.720861031978093
$ws.Range("D13").Value = 8.322079568072784
$ws.Range("E13").Value = 13.29097376639992
$ws.Range("F13").Value = 41.02195977649676
$ws.Range("G13").Value = 3.618420080178632
$ws.Range("J13").Value = 9.988978231202863
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 31.10162943268654
$ws.Range("C14").Value = 5.687236925558879
$ws.Range("D14").Value = 8.32410381398088
$ws.Range("E14").Value = 13.28420635488916
$ws.Range("F14").Value = 40.90085237831989
$ws.Range("G14").Value = 3.619181755378346
$ws.Range("J14").Value = 9.987127206179927
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 31.01418544175401
$ws.Range("C15").Value = 5.666573847298719
$ws.Range("D15").Value = 8.325361550290696
$ws.Range("E15").Value = 13.28011634602458
$ws.Range("F15").Value = 40.82672936358039
$ws.Range("G15").Value = 3.619650890739402
$ws.Range("J15").Value = 9.986031596312477
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 30.96067881572731
$ws.Range("C16").Value = 5.547162608248234
$ws.Range("D16").Value = 8.332842487610495
$ws.Range("E16").Value = 13.25752713912069
$ws.Range("F16").Value = 40.402982238279
$ws.Range("G16").Value = 3.62237907425135
$ws.Range("J16").Value = 9.980338910971588
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 30.65498779838054
$ws.Range("C17").Value = 5.473052887979496
$ws.Range("D17").Value = 8.337675848305945
$ws.Range("E17").Value = 13.24442954544489
$ws.Range("F17").Value = 40.14407159247448
$ws.Range("G17").Value = 3.624088314602401
$ws.Range("J17").Value = 9.977369701394315
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 30.46838297396374
$ws.Range("C18").Value = 5.430128156078403
$ws.Range("D18").Value = 8.340545605306687
$ws.Range("E18").Value = 13.23717647786342
$ws.Range("F18").Value = 39.99558975782951
$ws.Range("G18").Value = 3.625084546993187
$ws.Range("J18").Value = 9.975854225933524
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 30.36143056204226
$ws.Range("C19").Value = 5.415544782208928
$ws.Range("D19").Value = 8.341532672509526
$ws.Range("E19").Value = 13.23476894553345
$ws.Range("F19").Value = 39.94539640463729
$ws.Range("G19").Value = 3.625424112149621
$ws.Range("J19").Value = 9.975374153824932
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 30.32528672754184
$ws.Range("C20").Value = 5.480973273934539
$ws.Range("D20").Value = 8.337152042595166
$ws.Range("E20").Value = 13.24579481525483
$ws.Range("F20").Value = 40.17158907014755
$ws.Range("G20").Value = 3.623905005908867
$ws.Range("J20").Value = 9.977665874913908
$ws.Range("N20").Value = 20.96544799483449
$ws.Range("O20").Value = 30.48820912084517
$ws.Range("C21").Value = 5.697127250980179
$ws.Range("D21").Value = 8.323505528986626
$ws.Range("E21").Value = 13.28618260483261
$ws.Range("F21").Value = 40.93641259953594
$ws.Range("G21").Value = 3.618957493959642
$ws.Range("J21").Value = 9.987662937952059
$ws.Range("N21").Value = 21.70751365554064
$ws.Range("O21").Value = 31.03985853224028
$ws.Range("C22").Value = 5.835422455774041
$ws.Range("D22").Value = 8.315384781008332
$ws.Range("E22").Value = 13.31505285444418
$ws.Range("F22").Value = 41.43904107124759
$ws.Range("G22").Value = 3.615839929080495
$ws.Range("J22").Value = 9.995906512116415
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 31.40297030428751
$ws.Range("C23").Value = 5.761897942792122
$ws.Range("D23").Value = 8.319646008059562
$ws.Range("E23").Value = 13.29941899025567
$ws.Range("F23").Value = 41.17057863238634
$ws.Range("G23").Value = 3.617493270602819
$ws.Range("J23").Value = 9.991350479634942
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 31.20897271189365
$ws.Range("C24").Value = 5.477393454533012
$ws.Range("D24").Value = 8.337388571721
$ws.Range("E24").Value = 13.24517671402121
$ws.Range("F24").Value = 40.1591472555917
$ws.Range("G24").Value = 3.623987837499461
$ws.Range("J24").Value = 9.977531378164608
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 30.47924468572222
$ws.Range("C25").Value = 5.159325451038389
$ws.Range("D25").Value = 8.359958167408955
$ws.Range("E25").Value = 13.19739254666056
$ws.Range("F25").Value = 39.08552773273382
$ws.Range("G25").Value = 3.631495555123405
$ws.Range("J25").Value = 9.969978108515274
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 29.70702748933676
